# Commit: "Revised Prevention or Preventative Health to Prevention"
#
# The chiropractic_care table has a column header in cell C1 that read
# "Preventative" (paired with "At Risk", "Sick Role", "Self Care" as the
# motivation categories). It is revised to read "Prevention" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chiropractic_care.xlsx")

$ws.Range("C1").Value = "Prevention"
